# "mona minta yang lulus jadi layak"
# Flip the "label" column (G) on Sheet1:
#   - rows that were 1 (lulus) -> 0
#   - rows that were 0 (tidak lulus) -> 1 (layak)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose G value flips from 1 -> 0
$rowsToZero = @(4, 25, 29, 34, 37)

# Rows whose G value flips from 0 -> 1
$rowsToOne = @(
    47, 48, 49, 50, 51, 52, 53,
    55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67,
    69, 70, 71, 72,
    74, 75, 76, 77, 78,
    80, 81
)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 7).Value = 0
}

foreach ($r in $rowsToOne) {
    $ws.Cells.Item($r, 7).Value = 1
}

# Update the active selection on the sheet to match the resulting state
$ws.Range("G2:G81").Select()
